$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old sample rows (2-4) entirely and the hyperlinks attached to them ---
$ws.Hyperlinks.Delete()
$ws.Range("A2:AO4").Clear()

# --- New data: ten "Change" rows that all report the same outcome ---
$productIds = @(
    "TEST - Dummy 01",
    "TEST - Dummy 10",
    "TEST - Dummy 04",
    "TEST - Dummy 05",
    "TEST - Dummy 06",
    "TEST - Dummy 07",
    "TEST - Dummy 08",
    "TEST - Dummy 09",
    "TEST - Dummy 02",
    "TEST - Dummy 15 "
)

$row = 2
foreach ($productId in $productIds) {
    $ws.Cells.Item($row, 1).Value = "Change"
    $row = $row + 1
}

$row = 2
foreach ($productId in $productIds) {
    $ws.Cells.Item($row, 4).Value = $productId
    $row = $row + 1
}

$row = 2
foreach ($productId in $productIds) {
    $ws.Cells.Item($row, 15).Value = "y"
    $row = $row + 1
}

$row = 2
foreach ($productId in $productIds) {
    $ws.Cells.Item($row, 16).Value = 100
    $row = $row + 1
}

$row = 2
foreach ($productId in $productIds) {
    $ws.Cells.Item($row, 19).Value = "mis@contoso.com;supv@contoso.com"
    $row = $row + 1
}

# --- Hyperlink each Notify Emails cell to the mailto target ---
$row = 2
foreach ($productId in $productIds) {
    $cell = $ws.Cells.Item($row, 19)
    $ws.Hyperlinks.Add($cell, "mailto:mis@contoso.com;supv@contoso.com")
    $cell.Style = "Hyperlink"
    $row = $row + 1
}

$ws.Range("E11").Select()
